$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.306.51'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.873.54'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7085'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07801'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3110'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.884.52'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7183'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.16'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008388'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.136'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.75%  '
$ws.Range('D18').Value = '29.316.14'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.46%  '
$ws.Range('D20').Value = '2.130.25'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.032'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.504'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.418'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.333'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.243'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05344'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7525'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = '1.241.22'
$ws.Range('E39').Value = '  +6.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.733'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.526'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8951'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '109.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000127'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.52%  '
$ws.Range('D47').Value = '2.021.18'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5199'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.791'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.449'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4345'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.79%  '
